$d = $word.ActiveDocument

# 1) Title paragraph: merge "Research objectives and use of obtained " + "data"
#    (the gramStart/gramEnd proofErr marks around "data" go away with the edit)
$d.Content.Find.Execute(
    "Research objectives and use of obtained data",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Research objectives and use of obtained data",
    2) | Out-Null

# 2) "run-of-river" -> "Run-of-River (RoR)"
$d.Content.Find.Execute(
    "run-of-river hydroelectric systems",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Run-of-River (RoR) hydroelectric systems",
    2) | Out-Null

# 3) Remove the sentence "The energy yield will be calculated using the same
#    formula used by the hydro operators during their estimation. " that sat
#    between "...for a given time period. " and "For both ".
$d.Content.Find.Execute(
    "for a given time period. The energy yield will be calculated using the same formula used by the hydro operators during their estimation. For both ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "for a given time period. For both ",
    2) | Out-Null

# 4) Merge "available" + " " run split into a single "available " run, and
#    append the new closing sentence about peaking RoR projects after
#    "measured values."
$d.Content.Find.Execute(
    "available measured values.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "available measured values. For the peaking RoR projects, decision making scenarios is studied as they allow short term water storage option. ",
    2) | Out-Null
